# Slide 2, "Content Placeholder 2" (Shapes.Item(2)):
#   - Merge the "Areas of "/"Experties"/":" runs into a single corrected
#     run "Areas of Expertise:".
#   - Merge the "Payments, "/"Fintechs"/", Banking and ..." runs into a
#     single run with "Fintechs" corrected to "FinTech".
#
# Setting .Text directly to the corrected wording would just patch the
# differing characters in-place and keep the original multi-run split
# (PowerPoint preserves formatting boundaries for the common prefix).
# Assigning an unrelated placeholder string first collapses the
# paragraph down to one run (taking the formatting of the original first
# run), and then the follow-up assignment to the real text no longer
# shares a prefix with multiple old runs, so it stays a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

$paraHeading = $tr.Paragraphs(7)
$paraHeading.Text = "placeholder"
$paraHeading.Text = "Areas of Expertise:"

$paraDetail = $tr.Paragraphs(8)
$paraDetail.Text = "placeholder"
$paraDetail.Text = "Payments, FinTech, Banking and Finance, Branchless Banking, Software Engineering and Architecture, Application Desing, UI/UX Invision, Business Intelligence."
